$d = $word.ActiveDocument

function Replace-AllText($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Simple text replacements (appear once or twice, all occurrences map to the same new value)
Replace-AllText "2344567" "2487097"
Replace-AllText "01.06.2024" "03.05.2024"
Replace-AllText "One Big Chan Co. Ltd" "Акционерное общество «Арктические морские инженерно-геологические экспедиции»"
Replace-AllText "Шанхай, Китай" "Светлый"
Replace-AllText "инженер-инспектор Добрусев П. Ю." "старший инженер-инспектор Козлов С. В."
Replace-AllText "Доверенности № 12345 от 17.03.2024" "Доверенности № 122 от 31.01.2024"
Replace-AllText "капитан Иванов Р. И." "Капитан Бахтин Ю. Г."
Replace-AllText "911287" "940330"
Replace-AllText "Ежегодное освидетельствование" " освидетельствование"
Replace-AllText "1 000 000,00 p. (один миллион рублей 00 копеек)" "10 000,00 p. (десять тысяч рублей 00 копеек)"
Replace-AllText "Р. И. Иванов" "Ю. Г. Бахтин"
Replace-AllText "П. Ю. Добрусев" "С. В. Козлов"

# Vessel name/number: contains straight quotes - set Range.Text directly to avoid
# Word's "smart quotes" AutoCorrect turning " into curly quotes during Find/Replace.
$rng = $d.Content
if ($rng.Find.Execute('"ВЛАДИВОСТОК" / "VLADIVOSTOK" ', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Text = '"ВОЛГА" '
}

# The empty run right before "(documents issued by RS ...)" needs new text inserted
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -like "*documents issued*") {
            $p.Range.InsertBefore("Свидетельство ф. 6.5.30 №№ 24.42.03.00234.121 - 24.42.03.00236.121 от --")
            break
        }
    }
}
